$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above the existing row 249 ("Vega Central Mapocho de
# Santiago" / "Rabanito" records). This pushes the previous rows 249:310
# down to 250:311, matching the target dimension A1:R311.
$ws.Rows.Item(249).Insert()

# Populate the newly inserted row 249 with a new weekly observation. The
# "static" columns (A,B,C,E,F,G,H,I,N,Q,R) mirror every other row in this
# table; K,L,M,O,P carry over the values that used to be on row 249 before
# the insert (now row 250), while D (fecha) and J (volumen) are the new
# data points.
$ws.Range("A249").Value = 9
$ws.Range("B249").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C249").Value = "Metropolitana"
$ws.Range("D249").Value = 44798
$ws.Range("E249").Value = 13
$ws.Range("F249").Value = 300000001
$ws.Range("G249").Value = "Rabanito"
$ws.Range("H249").Value = "Sin especificar"
$ws.Range("I249").Value = "Primera"
$ws.Range("J249").Value = 7900
$ws.Range("K249").Value = 3000
$ws.Range("L249").Value = 3000
$ws.Range("M249").Value = 3000
$ws.Range("N249").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O249").Value = "Provincia de Chacabuco"
$ws.Range("P249").Value = 30
$ws.Range("Q249").Value = 100
$ws.Range("R249").Value = "Hortaliza"
